# Add the new "2022-Q3" sheet, positioned right after "总计" and before "2022-Q2".
# Existing sheets (2022-Q2 .. 2020-Q4) automatically shift right; their sheet
# content is otherwise unchanged.
$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

function Write-FundRow {
    param($sheet, $row, $values)
    for ($c = 1; $c -le $values.Length; $c++) {
        $cell = $values[$c - 1]
        if ($null -eq $cell) { continue }
        $kind = $cell.Substring(0, 2)
        $val = $cell.Substring(2)
        $rng = $sheet.Cells.Item($row, $c)
        if ($kind -eq "S:") {
            # Force text format first so numeric-looking strings (fund codes
            # with leading zeros, "95.85", ...) are written as literal text
            # instead of silently turning into numbers; then drop the
            # left-over "@" number-format style so the cell format matches
            # the plain (unstyled) text cells used throughout the workbook.
            $rng.NumberFormat = "@"
            $rng.Value = $val
            $rng.ClearFormats()
        } else {
            $rng.Value = [double]$val
        }
    }
}

$q3Data = @(
  ,($null, "S:基金代码", "S:基金名称", "S:基金规模", "S:股票总仓位", "S:仓位占比", "S:持有市值(亿元)", "S:仓位排名")
  ,("N:0", "S:090018", "S:大成新锐产业混合", "S:95.85", "S:83.49", "S:7.72", "S:7.3996", "N:3")
  ,("N:1", "S:001300", "S:大成睿景灵活配置混合A", "S:33.29", "S:90.96", "S:8.35", "S:2.7797", "N:3")
  ,("N:2", "S:013435", "S:大成景气精选六个月持有混合A", "S:32.65", "S:88.02", "S:6.49", "S:2.1190", "N:3")
  ,("N:3", "S:001301", "S:大成睿景灵活配置混合C", "S:23.86", "S:90.96", "S:8.35", "S:1.9923", "N:3")
  ,("N:4", "S:002258", "S:大成国企改革灵活配置混合", "S:17.20", "S:90.32", "S:8.29", "S:1.4259", "N:3")
  ,("N:5", "S:014224", "S:大成聚优成长混合A", "S:16.94", "S:86.11", "S:6.94", "S:1.1756", "N:3")
  ,("N:6", "S:010826", "S:大成产业趋势混合A", "S:11.59", "S:90.99", "S:8.31", "S:0.9631", "N:3")
  ,("N:7", "S:012519", "S:大成核心趋势混合A", "S:9.97", "S:88.72", "S:8.75", "S:0.8724", "N:3")
  ,("N:8", "S:010827", "S:大成产业趋势混合C", "S:6.29", "S:90.99", "S:8.31", "S:0.5227", "N:3")
  ,("N:9", "S:008934", "S:大成科技消费股票A", "S:10.12", "S:84.64", "S:4.11", "S:0.4159", "N:6")
  ,("N:10", "S:013436", "S:大成景气精选六个月持有混合C", "S:5.65", "S:88.02", "S:6.49", "S:0.3667", "N:3")
  ,("N:11", "S:012184", "S:大成创新趋势混合A", "S:7.68", "S:63.97", "S:4.20", "S:0.3226", "N:3")
  ,("N:12", "S:014225", "S:大成聚优成长混合C", "S:3.36", "S:86.11", "S:6.94", "S:0.2332", "N:3")
  ,("N:13", "S:160918", "S:大成中小盘混合（LOF）A", "S:5.28", "S:63.53", "S:4.08", "S:0.2154", "N:3")
  ,("N:14", "S:012520", "S:大成核心趋势混合C", "S:2.45", "S:88.72", "S:8.75", "S:0.2144", "N:3")
  ,("N:15", "S:014185", "S:招商专精特新股票A", "S:3.62", "S:81.87", "S:4.80", "S:0.1738", "N:9")
  ,("N:16", "S:008935", "S:大成科技消费股票C", "S:3.90", "S:84.64", "S:4.11", "S:0.1603", "N:6")
  ,("N:17", "S:014186", "S:招商专精特新股票C", "S:2.63", "S:81.87", "S:4.80", "S:0.1262", "N:9")
  ,("N:18", "S:008274", "S:大成行业先锋混合A", "S:2.60", "S:67.57", "S:4.21", "S:0.1095", "N:3")
  ,("N:19", "S:217013", "S:招商中小盘精选混合", "S:2.52", "S:80.73", "S:4.26", "S:0.1074", "N:8")
  ,("N:20", "S:002945", "S:大成盛世精选灵活配置混合", "S:1.25", "S:63.59", "S:4.19", "S:0.0524", "N:3")
  ,("N:21", "S:008275", "S:大成行业先锋混合C", "S:0.49", "S:67.57", "S:4.21", "S:0.0206", "N:3")
  ,("N:22", "S:001531", "S:招商安益灵活配置混合", "S:0.57", "S:70.35", "S:3.45", "S:0.0197", "N:10")
  ,("N:23", "S:540007", "S:汇丰晋信中小盘股票", "S:0.56", "S:92.05", "S:3.37", "S:0.0189", "N:8")
  ,("N:24", "S:012185", "S:大成创新趋势混合C", "S:0.15", "S:63.97", "S:4.20", "S:0.0063", "N:3")
  ,("N:25", "S:011159", "S:大成中小盘混合（LOF）C", "S:0.01", "S:63.53", "S:4.08", "S:0.0004", "N:3")
)

$r = 1
foreach ($row in $q3Data) {
    Write-FundRow $q3Sheet $r $row
    $r = $r + 1
}

# Column A (row index) and the header row (row 1) use the bold / centered /
# thin-bordered "index" style in every fund sheet - reproduce that look.
$q3HeaderAndIndexCells = @("B1","C1","D1","E1","F1","G1","H1")
for ($i = 2; $i -le 27; $i++) {
    $q3HeaderAndIndexCells += "A$i"
}
foreach ($ref in $q3HeaderAndIndexCells) {
    $rng = $q3Sheet.Range($ref)
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# Insert the new "2022-Q3" summary row into "总计" (row 2), pushing the
# existing quarters down by one row.
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Match column A's existing bold/centered/bordered "index" style by copying
# the format from the row below (already-correct existing data row) before
# writing the new index value.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 26
$totalSheet.Range("D2").Value = 21.81
